$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.784564733505249
$ws.Range("B1").Value = 2.263836622238159
$ws.Range("C1").Value = 2.256114482879639
$ws.Range("D1").Value = 1.998031735420227
$ws.Range("E1").Value = 1.357514381408691
